$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M2").Value = 0.1811433333333334
$ws.Range("N2").Value = 0.5434300000000001
$ws.Range("O2").Value = 0.0111261749556462
$ws.Range("P2").Value = 0.01112617495564619
$ws.Range("Q2").Value = 0.07648940279000001
$ws.Range("R2").Value = 0.6884046251100001
$ws.Range("S2").Value = 0.0111261749556462
$ws.Range("T2").Value = 0.01112617495564619

$ws.Range("O3").Value = 0.8246098959508241
$ws.Range("P3").Value = 0.8246098959508241
$ws.Range("Q3").Value = 5.668966983481999
$ws.Range("R3").Value = 51.020702851338
$ws.Range("S3").Value = 0.8246098959508241
$ws.Range("T3").Value = 0.8246098959508241

$ws.Range("M4").Value = 2.659118666666667
$ws.Range("N4").Value = 7.977356
$ws.Range("O4").Value = 0.1633282272592126
$ws.Range("P4").Value = 0.1633282272592126
$ws.Range("Q4").Value = 1.122836789068
$ws.Range("R4").Value = 10.105531101612
$ws.Range("S4").Value = 0.1633282272592126
$ws.Range("T4").Value = 0.1633282272592126

$ws.Range("M5").Value = 0.015234
$ws.Range("N5").Value = 0.045702
$ws.Range("O5").Value = 0.0009357018343171013
$ws.Range("P5").Value = 0.0009357018343171013
$ws.Range("Q5").Value = 0.006432693606
$ws.Range("R5").Value = 0.057894242454
$ws.Range("S5").Value = 0.0009357018343171013
$ws.Range("T5").Value = 0.0009357018343171013
